$d = $word.ActiveDocument

# Make clean edits (no tracked-change markup), restoring the document's
# original TrackRevisions setting afterwards so we don't leave an
# unrelated change behind in word/settings.xml.
$origTrack = $d.TrackRevisions
$d.TrackRevisions = $false

# Update the three estimate figures in the table.
$d.Content.Find.Execute("42,037", $true, $false, $false, $false, $false, $true, 1, $false, "44,183", 2) | Out-Null
$d.Content.Find.Execute("55,500", $true, $false, $false, $false, $false, $true, 1, $false, "56,638", 2) | Out-Null
$d.Content.Find.Execute("82,620", $true, $false, $false, $false, $false, $true, 1, $false, "74,344", 2) | Out-Null

# Re-write the "2020" cell in place; this drops the stale
# <w:lastRenderedPageBreak/> marker on that run (Word regenerates these
# only on the next real pagination, so a clean replace removes it),
# matching the diff exactly.
$d.Content.Find.Execute("2020", $true, $false, $false, $false, $false, $true, 1, $false, "2020", 2) | Out-Null

$d.TrackRevisions = $origTrack
